# Add architecture and testbench of Top Unit
# Appends a new time-record row (row 10) for 29.11.2019 / "Top Unit",
# mirroring the existing rows (date, start time, end time, duration
# formula, category, description).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date of the new entry
$ws.Range("A10").Value = "29.11.2019"

# Start / end time (10:00 - 11:30), using the same time format as the
# other "Architecture, Testbench" rows (B6:C9 use the h:mm style).
$ws.Range("B10").Value = 0.41666666666666669
$ws.Range("B10").NumberFormat = "h:mm"

$ws.Range("C10").Value = 0.47916666666666669
$ws.Range("C10").NumberFormat = "h:mm"

# Duration = end - start, formatted like the rest of column D.
$ws.Range("D10").Formula = "=C10-B10"
$ws.Range("D10").NumberFormat = "[$]hh:mm;@"

# Category / unit worked on, and description of the work done.
$ws.Range("E10").Value = "Top Unit"
$ws.Range("F10").Value = "Architecture, Testbench"

# Move the active selection onto the newly added cell, like a user
# would after typing the last value of the new row.
$ws.Range("F10").Select()
